$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 72.14286
$ws.Range("I5").Value = 77.5
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 77.5
$ws.Range("L5").Value = 40
$ws.Range("M5").Value = 37.5
$ws.Range("N5").Value = -270

$ws.Range("H12").Value = 117
$ws.Range("I12").Value = 117
$ws.Range("K12").Value = 117
$ws.Range("M12").Value = 53

$ws.Range("H28").Value = 856.86664
$ws.Range("I28").Value = 612.8333
$ws.Range("K28").Value = 612.8333
$ws.Range("M28").Value = -127.8333

$ws.Range("H33").Value = 299.42856
$ws.Range("J33").Value = 644.5
$ws.Range("L33").Value = 644.5
$ws.Range("N33").Value = -1102.5

$ws.Range("H70").Value = 8864.777
$ws.Range("J70").Value = 9598
$ws.Range("L70").Value = 28794
$ws.Range("N70").Value = -29334

$ws.Range("H73").Value = 8864.777
$ws.Range("J73").Value = 9598
$ws.Range("L73").Value = 28794
$ws.Range("N73").Value = -30666

$ws.Range("H116").Value = 24454.834
$ws.Range("I116").Value = 4530.2856
$ws.Range("K116").Value = 4530.2856
$ws.Range("M116").Value = -1088.2856

$ws.Range("H127").Value = 3544.4092
$ws.Range("I127").Value = 3636.158
$ws.Range("J127").Value = 2963.3333
$ws.Range("K127").Value = 10908.474
$ws.Range("L127").Value = 8889.999899999999
$ws.Range("M127").Value = -5948.474
$ws.Range("N127").Value = -18809.9999

$ws.Range("H132").Value = 11439.091
$ws.Range("I132").Value = 11602.857
$ws.Range("K132").Value = 34808.571
$ws.Range("M132").Value = -32278.571

$ws.Range("H137").Value = 10343.95
$ws.Range("I137").Value = 2249.3
$ws.Range("J137").Value = 18438.6
$ws.Range("K137").Value = 6747.900000000001
$ws.Range("L137").Value = 55315.8
$ws.Range("M137").Value = -4197.900000000001
$ws.Range("N137").Value = -60415.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 131387.77
$ws.Range("I32").Value = 134281
$ws.Range("J32").Value = 19998
$ws.Range("K32").Value = 134281
$ws.Range("L32").Value = 19998
$ws.Range("M32").Value = -133994
$ws.Range("N32").Value = -20572

$ws.Range("H88").Value = 1685.2858
$ws.Range("J88").Value = 1853.3158
$ws.Range("L88").Value = 1853.3158
$ws.Range("N88").Value = -2665.3158

$ws.Range("H91").Value = 1685.2858
$ws.Range("J91").Value = 1853.3158
$ws.Range("L91").Value = 1853.3158
$ws.Range("N91").Value = -4661.3158

$ws.Range("H102").Value = 9517.412
$ws.Range("I102").Value = 9517.412
$ws.Range("K102").Value = 9517.412
$ws.Range("M102").Value = -7895.412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 59796
$ws.Range("J40").Value = 59796
$ws.Range("L40").Value = 59796
$ws.Range("N40").Value = -60326

$ws.Range("H96").Value = 40793.668
$ws.Range("I96").Value = 31203
$ws.Range("J96").Value = 59975
$ws.Range("K96").Value = 31203
$ws.Range("L96").Value = 59975
$ws.Range("M96").Value = -28457
$ws.Range("N96").Value = -65467

$ws.Range("H134").Value = 10873.6
$ws.Range("I134").Value = 6420.6665
$ws.Range("K134").Value = 19261.9995
$ws.Range("M134").Value = -16726.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 55.625
$ws.Range("I7").Value = 40.833332
$ws.Range("K7").Value = 40.833332
$ws.Range("M7").Value = 72.166668

$ws.Range("H10").Value = 760.4545000000001
$ws.Range("I10").Value = 818.3333
$ws.Range("K10").Value = 818.3333
$ws.Range("M10").Value = -679.3333

$ws.Range("H22").Value = 1555.3889
$ws.Range("I22").Value = 636.8182
$ws.Range("K22").Value = 636.8182
$ws.Range("M22").Value = -286.8182

$ws.Range("H31").Value = 3532.2563
$ws.Range("I31").Value = 4063.818
$ws.Range("J31").Value = 2844.353
$ws.Range("K31").Value = 4063.818
$ws.Range("L31").Value = 2844.353
$ws.Range("M31").Value = -3768.818
$ws.Range("N31").Value = -3434.353

$ws.Range("H34").Value = 3532.2563
$ws.Range("I34").Value = 4063.818
$ws.Range("J34").Value = 2844.353
$ws.Range("K34").Value = 4063.818
$ws.Range("L34").Value = 2844.353
$ws.Range("M34").Value = -3861.818
$ws.Range("N34").Value = -3248.353

$ws.Range("H58").Value = 11268.875
$ws.Range("I58").Value = 4928.4546
$ws.Range("K58").Value = 4928.4546
$ws.Range("M58").Value = -4725.4546

$ws.Range("H107").Value = 357.42856
$ws.Range("J107").Value = 517.8
$ws.Range("L107").Value = 517.8
$ws.Range("N107").Value = -4357.8

$ws.Range("H132").Value = 3471.3845
$ws.Range("I132").Value = 3471.3845
$ws.Range("K132").Value = 10414.1535
$ws.Range("M132").Value = -7884.1535

$ws.Range("H136").Value = 11268.875
$ws.Range("I136").Value = 4928.4546
$ws.Range("K136").Value = 14785.3638
$ws.Range("M136").Value = -12235.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 9856.429
$ws.Range("J127").Value = 9856.429
$ws.Range("L127").Value = 29569.287
$ws.Range("N127").Value = -39489.287

$ws.Range("H131").Value = 5430
$ws.Range("J131").Value = 8625.625
$ws.Range("L131").Value = 25876.875
$ws.Range("N131").Value = -35956.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 29999.5
$ws.Range("J10").Value = 29999.5
$ws.Range("L10").Value = 29999.5
$ws.Range("N10").Value = -30337.5

$ws.Range("H11").Value = 5784900
$ws.Range("J11").Value = 2705000.2
$ws.Range("L11").Value = 2705000.2
$ws.Range("N11").Value = -2705278.2

$ws.Range("H103").Value = 75000
$ws.Range("J103").Value = 75000
$ws.Range("L103").Value = 75000
$ws.Range("N103").Value = -77344

$ws.Range("H113").Value = 3017.4814
$ws.Range("I113").Value = 2261.9048
$ws.Range("K113").Value = 2261.9048
$ws.Range("M113").Value = -91.9047999999998

$ws.Range("H122").Value = 5111.6665
$ws.Range("I122").Value = 1001.2
$ws.Range("J122").Value = 10249.75
$ws.Range("K122").Value = 3003.6
$ws.Range("L122").Value = 30749.25
$ws.Range("M122").Value = -553.6000000000004
$ws.Range("N122").Value = -35649.25

$ws.Range("H123").Value = 49999
$ws.Range("J123").Value = 49999
$ws.Range("L123").Value = 49999
$ws.Range("N123").Value = -54899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9500.842000000001
$ws.Range("J136").Value = 12505.182
$ws.Range("L136").Value = 37515.546
$ws.Range("N136").Value = -42615.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4887

$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 12500
$ws.Range("K62").Value = 12500
$ws.Range("M62").Value = -11876

$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 12500
$ws.Range("K65").Value = 62500
$ws.Range("M65").Value = -59380

$ws.Range("H132").Value = 2804.087
$ws.Range("J132").Value = 3660
$ws.Range("L132").Value = 10980
$ws.Range("N132").Value = -16040

$ws.Range("H136").Value = 870.64703
$ws.Range("I136").Value = 877.46155
$ws.Range("K136").Value = 2632.38465
$ws.Range("M136").Value = -82.38464999999997
